$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("Y1").Copy()
$ws1.Range("Z1").PasteSpecial(-4122)
$ws1.Range("Z1").WrapText = $false
$ws1.Range("Z1").Value = "Loại cán sự"
Write-Output "done"
